$d = $word.ActiveDocument

$replacements = @(
    @{old = "2024-01-30 Tuesday"; new = "2024-01-31 Wednesday"},
    @{old = "86×32=2752"; new = "45×54=2430"},
    @{old = "93×76=7068"; new = "64×18=1152"},
    @{old = "21×47=987";  new = "57×98=5586"},
    @{old = "11×75=825";  new = "33×44=1452"},
    @{old = "63×29=1827"; new = "59×57=3363"},
    @{old = "51×80=4080"; new = "40×24=960"},
    @{old = "22×21=462";  new = "50×37=1850"},
    @{old = "45×39=1755"; new = "88×53=4664"},
    @{old = "88×90=7920"; new = "28×92=2576"},
    @{old = "34×72=2448"; new = "60×82=4920"},
    @{old = "25×25=625";  new = "12×71=852"},
    @{old = "58×28=1624"; new = "95×54=5130"},
    @{old = "71×50=3550"; new = "25×69=1725"},
    @{old = "57×82=4674"; new = "75×74=5550"},
    @{old = "68×99=6732"; new = "19×18=342"},
    @{old = "76×52=3952"; new = "73×27=1971"},
    @{old = "64×63=4032"; new = "64×19=1216"},
    @{old = "11×42=462";  new = "80×49=3920"},
    @{old = "98×74=7252"; new = "25×77=1925"},
    @{old = "77×95=7315"; new = "63×61=3843"},
    @{old = "22×36=792";  new = "74×84=6216"},
    @{old = "71×29=2059"; new = "68×69=4692"},
    @{old = "44×23=1012"; new = "24×85=2040"},
    @{old = "46×13=598";  new = "74×62=4588"},
    @{old = "80×83=6640"; new = "11×16=176"}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $true, $false, $false, $false, $true, 1, $false, $r.new, 2)
}

$d.Save()
